$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows ---
$ws.Cells.Item(56, 17).Value2 = 0   # Q56: 1 -> 0
$ws.Cells.Item(68, 17).Value2 = 0   # Q68: 1 -> 0
$ws.Cells.Item(69, 17).Value2 = 0   # Q69: 2 -> 0
$ws.Cells.Item(1485, 15).Value2 = 1 # O1485: 0 -> 1
$ws.Cells.Item(1487, 18).Value2 = 0 # R1487: "" -> 0 (numeric)
$ws.Cells.Item(1488, 18).Value2 = 0 # R1488: "" -> 0 (numeric)

# --- Append new rows 1489-1506 ---
# Row 1489
$ws.Cells.Item(1489, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1489, 1).Value2 = 45474
$ws.Cells.Item(1489, 2).Value2 = 2461.050048828125
$ws.Cells.Item(1489, 3).Value2 = 2554
$ws.Cells.Item(1489, 4).Value2 = 2450.10009765625
$ws.Cells.Item(1489, 5).Value2 = 2547
$ws.Cells.Item(1489, 6).Value2 = 2547
$ws.Cells.Item(1489, 7).Value2 = 9755776
$ws.Cells.Item(1489, 8).Value2 = 2024
$ws.Cells.Item(1489, 9).Value2 = 7
$ws.Cells.Item(1489, 10).Value2 = 1
$ws.Cells.Item(1489, 11).Value2 = 0
$ws.Cells.Item(1489, 12).Value2 = 0
$ws.Cells.Item(1489, 13).Value2 = 0
$ws.Cells.Item(1489, 14).Value2 = 27
$ws.Cells.Item(1489, 15).Value2 = 0
$ws.Cells.Item(1489, 16).Value2 = 0
$ws.Cells.Item(1489, 17).Value2 = 0

# Row 1490
$ws.Cells.Item(1490, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1490, 1).Value2 = 45481
$ws.Cells.Item(1490, 2).Value2 = 2560
$ws.Cells.Item(1490, 3).Value2 = 2639.5
$ws.Cells.Item(1490, 4).Value2 = 2536.25
$ws.Cells.Item(1490, 5).Value2 = 2622.25
$ws.Cells.Item(1490, 6).Value2 = 2622.25
$ws.Cells.Item(1490, 7).Value2 = 9407687
$ws.Cells.Item(1490, 8).Value2 = 2024
$ws.Cells.Item(1490, 9).Value2 = 7
$ws.Cells.Item(1490, 10).Value2 = 8
$ws.Cells.Item(1490, 11).Value2 = 0
$ws.Cells.Item(1490, 12).Value2 = 0
$ws.Cells.Item(1490, 13).Value2 = 0
$ws.Cells.Item(1490, 14).Value2 = 28
$ws.Cells.Item(1490, 15).Value2 = 0
$ws.Cells.Item(1490, 16).Value2 = 0
$ws.Cells.Item(1490, 17).Value2 = 0

# Row 1491
$ws.Cells.Item(1491, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1491, 1).Value2 = 45488
$ws.Cells.Item(1491, 2).Value2 = 2628
$ws.Cells.Item(1491, 3).Value2 = 2751.199951171875
$ws.Cells.Item(1491, 4).Value2 = 2613.64990234375
$ws.Cells.Item(1491, 5).Value2 = 2727
$ws.Cells.Item(1491, 6).Value2 = 2727
$ws.Cells.Item(1491, 7).Value2 = 11456728
$ws.Cells.Item(1491, 8).Value2 = 2024
$ws.Cells.Item(1491, 9).Value2 = 7
$ws.Cells.Item(1491, 10).Value2 = 15
$ws.Cells.Item(1491, 11).Value2 = 0
$ws.Cells.Item(1491, 12).Value2 = 0
$ws.Cells.Item(1491, 13).Value2 = 0
$ws.Cells.Item(1491, 14).Value2 = 29
$ws.Cells.Item(1491, 15).Value2 = 0
$ws.Cells.Item(1491, 16).Value2 = 0
$ws.Cells.Item(1491, 17).Value2 = 2

# Row 1492
$ws.Cells.Item(1492, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1492, 1).Value2 = 45495
$ws.Cells.Item(1492, 2).Value2 = 2705.699951171875
$ws.Cells.Item(1492, 3).Value2 = 2811.300048828125
$ws.Cells.Item(1492, 4).Value2 = 2670
$ws.Cells.Item(1492, 5).Value2 = 2711.050048828125
$ws.Cells.Item(1492, 6).Value2 = 2711.050048828125
$ws.Cells.Item(1492, 7).Value2 = 13731794
$ws.Cells.Item(1492, 8).Value2 = 2024
$ws.Cells.Item(1492, 9).Value2 = 7
$ws.Cells.Item(1492, 10).Value2 = 22
$ws.Cells.Item(1492, 11).Value2 = 0
$ws.Cells.Item(1492, 12).Value2 = 0
$ws.Cells.Item(1492, 13).Value2 = 0
$ws.Cells.Item(1492, 14).Value2 = 30
$ws.Cells.Item(1492, 15).Value2 = 0
$ws.Cells.Item(1492, 16).Value2 = 0
$ws.Cells.Item(1492, 17).Value2 = 0

# Row 1493
$ws.Cells.Item(1493, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1493, 1).Value2 = 45502
$ws.Cells.Item(1493, 2).Value2 = 2711.199951171875
$ws.Cells.Item(1493, 3).Value2 = 2759.5
$ws.Cells.Item(1493, 4).Value2 = 2675.050048828125
$ws.Cells.Item(1493, 5).Value2 = 2692.550048828125
$ws.Cells.Item(1493, 6).Value2 = 2692.550048828125
$ws.Cells.Item(1493, 7).Value2 = 7893759
$ws.Cells.Item(1493, 8).Value2 = 2024
$ws.Cells.Item(1493, 9).Value2 = 7
$ws.Cells.Item(1493, 10).Value2 = 29
$ws.Cells.Item(1493, 11).Value2 = 0
$ws.Cells.Item(1493, 12).Value2 = 0
$ws.Cells.Item(1493, 13).Value2 = 0
$ws.Cells.Item(1493, 14).Value2 = 31
$ws.Cells.Item(1493, 15).Value2 = 0
$ws.Cells.Item(1493, 16).Value2 = 0
$ws.Cells.Item(1493, 17).Value2 = 0

# Row 1494
$ws.Cells.Item(1494, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1494, 1).Value2 = 45509
$ws.Cells.Item(1494, 2).Value2 = 2675.949951171875
$ws.Cells.Item(1494, 3).Value2 = 2781.85009765625
$ws.Cells.Item(1494, 4).Value2 = 2666.199951171875
$ws.Cells.Item(1494, 5).Value2 = 2747.199951171875
$ws.Cells.Item(1494, 6).Value2 = 2747.199951171875
$ws.Cells.Item(1494, 7).Value2 = 9876829
$ws.Cells.Item(1494, 8).Value2 = 2024
$ws.Cells.Item(1494, 9).Value2 = 8
$ws.Cells.Item(1494, 10).Value2 = 5
$ws.Cells.Item(1494, 11).Value2 = 0
$ws.Cells.Item(1494, 12).Value2 = 0
$ws.Cells.Item(1494, 13).Value2 = 0
$ws.Cells.Item(1494, 14).Value2 = 32
$ws.Cells.Item(1494, 15).Value2 = 0
$ws.Cells.Item(1494, 16).Value2 = 0
$ws.Cells.Item(1494, 17).Value2 = 0

# Row 1495
$ws.Cells.Item(1495, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1495, 1).Value2 = 45516
$ws.Cells.Item(1495, 2).Value2 = 2738
$ws.Cells.Item(1495, 3).Value2 = 2763.39990234375
$ws.Cells.Item(1495, 4).Value2 = 2706.60009765625
$ws.Cells.Item(1495, 5).Value2 = 2748.25
$ws.Cells.Item(1495, 6).Value2 = 2748.25
$ws.Cells.Item(1495, 7).Value2 = 5023336
$ws.Cells.Item(1495, 8).Value2 = 2024
$ws.Cells.Item(1495, 9).Value2 = 8
$ws.Cells.Item(1495, 10).Value2 = 12
$ws.Cells.Item(1495, 11).Value2 = 0
$ws.Cells.Item(1495, 12).Value2 = 0
$ws.Cells.Item(1495, 13).Value2 = 0
$ws.Cells.Item(1495, 14).Value2 = 33
$ws.Cells.Item(1495, 15).Value2 = 0
$ws.Cells.Item(1495, 16).Value2 = 0
$ws.Cells.Item(1495, 17).Value2 = 0

# Row 1496
$ws.Cells.Item(1496, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1496, 1).Value2 = 45523
$ws.Cells.Item(1496, 2).Value2 = 2750.050048828125
$ws.Cells.Item(1496, 3).Value2 = 2821
$ws.Cells.Item(1496, 4).Value2 = 2735.050048828125
$ws.Cells.Item(1496, 5).Value2 = 2815.60009765625
$ws.Cells.Item(1496, 6).Value2 = 2815.60009765625
$ws.Cells.Item(1496, 7).Value2 = 7310293
$ws.Cells.Item(1496, 8).Value2 = 2024
$ws.Cells.Item(1496, 9).Value2 = 8
$ws.Cells.Item(1496, 10).Value2 = 19
$ws.Cells.Item(1496, 11).Value2 = 0
$ws.Cells.Item(1496, 12).Value2 = 0
$ws.Cells.Item(1496, 13).Value2 = 0
$ws.Cells.Item(1496, 14).Value2 = 34
$ws.Cells.Item(1496, 15).Value2 = 0
$ws.Cells.Item(1496, 16).Value2 = 0
$ws.Cells.Item(1496, 17).Value2 = 0

# Row 1497
$ws.Cells.Item(1497, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1497, 1).Value2 = 45530
$ws.Cells.Item(1497, 2).Value2 = 2815.60009765625
$ws.Cells.Item(1497, 3).Value2 = 2834.949951171875
$ws.Cells.Item(1497, 4).Value2 = 2745
$ws.Cells.Item(1497, 5).Value2 = 2778
$ws.Cells.Item(1497, 6).Value2 = 2778
$ws.Cells.Item(1497, 7).Value2 = 10392651
$ws.Cells.Item(1497, 8).Value2 = 2024
$ws.Cells.Item(1497, 9).Value2 = 8
$ws.Cells.Item(1497, 10).Value2 = 26
$ws.Cells.Item(1497, 11).Value2 = 0
$ws.Cells.Item(1497, 12).Value2 = 0
$ws.Cells.Item(1497, 13).Value2 = 0
$ws.Cells.Item(1497, 14).Value2 = 35
$ws.Cells.Item(1497, 15).Value2 = 0
$ws.Cells.Item(1497, 16).Value2 = 0
$ws.Cells.Item(1497, 17).Value2 = 0

# Row 1498
$ws.Cells.Item(1498, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1498, 1).Value2 = 45537
$ws.Cells.Item(1498, 2).Value2 = 2794
$ws.Cells.Item(1498, 3).Value2 = 2866.449951171875
$ws.Cells.Item(1498, 4).Value2 = 2771.64990234375
$ws.Cells.Item(1498, 5).Value2 = 2838.949951171875
$ws.Cells.Item(1498, 6).Value2 = 2838.949951171875
$ws.Cells.Item(1498, 7).Value2 = 8238955
$ws.Cells.Item(1498, 8).Value2 = 2024
$ws.Cells.Item(1498, 9).Value2 = 9
$ws.Cells.Item(1498, 10).Value2 = 2
$ws.Cells.Item(1498, 11).Value2 = 0
$ws.Cells.Item(1498, 12).Value2 = 0
$ws.Cells.Item(1498, 13).Value2 = 0
$ws.Cells.Item(1498, 14).Value2 = 36
$ws.Cells.Item(1498, 15).Value2 = 0
$ws.Cells.Item(1498, 16).Value2 = 0
$ws.Cells.Item(1498, 17).Value2 = 0

# Row 1499
$ws.Cells.Item(1499, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1499, 1).Value2 = 45544
$ws.Cells.Item(1499, 2).Value2 = 2853.10009765625
$ws.Cells.Item(1499, 3).Value2 = 2963.39990234375
$ws.Cells.Item(1499, 4).Value2 = 2843.199951171875
$ws.Cells.Item(1499, 5).Value2 = 2932.949951171875
$ws.Cells.Item(1499, 6).Value2 = 2932.949951171875
$ws.Cells.Item(1499, 7).Value2 = 11149568
$ws.Cells.Item(1499, 8).Value2 = 2024
$ws.Cells.Item(1499, 9).Value2 = 9
$ws.Cells.Item(1499, 10).Value2 = 9
$ws.Cells.Item(1499, 11).Value2 = 0
$ws.Cells.Item(1499, 12).Value2 = 0
$ws.Cells.Item(1499, 13).Value2 = 0
$ws.Cells.Item(1499, 14).Value2 = 37
$ws.Cells.Item(1499, 15).Value2 = 0
$ws.Cells.Item(1499, 16).Value2 = 0
$ws.Cells.Item(1499, 17).Value2 = 0

# Row 1500
$ws.Cells.Item(1500, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1500, 1).Value2 = 45551
$ws.Cells.Item(1500, 2).Value2 = 2910
$ws.Cells.Item(1500, 3).Value2 = 2989.35009765625
$ws.Cells.Item(1500, 4).Value2 = 2807.39990234375
$ws.Cells.Item(1500, 5).Value2 = 2977.60009765625
$ws.Cells.Item(1500, 6).Value2 = 2977.60009765625
$ws.Cells.Item(1500, 7).Value2 = 8407156
$ws.Cells.Item(1500, 8).Value2 = 2024
$ws.Cells.Item(1500, 9).Value2 = 9
$ws.Cells.Item(1500, 10).Value2 = 16
$ws.Cells.Item(1500, 11).Value2 = 0
$ws.Cells.Item(1500, 12).Value2 = 0
$ws.Cells.Item(1500, 13).Value2 = 0
$ws.Cells.Item(1500, 14).Value2 = 38
$ws.Cells.Item(1500, 15).Value2 = 0
$ws.Cells.Item(1500, 16).Value2 = 0
$ws.Cells.Item(1500, 17).Value2 = 0

# Row 1501
$ws.Cells.Item(1501, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1501, 1).Value2 = 45558
$ws.Cells.Item(1501, 2).Value2 = 2980
$ws.Cells.Item(1501, 3).Value2 = 3035
$ws.Cells.Item(1501, 4).Value2 = 2905.10009765625
$ws.Cells.Item(1501, 5).Value2 = 2966.25
$ws.Cells.Item(1501, 6).Value2 = 2966.25
$ws.Cells.Item(1501, 7).Value2 = 9165102
$ws.Cells.Item(1501, 8).Value2 = 2024
$ws.Cells.Item(1501, 9).Value2 = 9
$ws.Cells.Item(1501, 10).Value2 = 23
$ws.Cells.Item(1501, 11).Value2 = 0
$ws.Cells.Item(1501, 12).Value2 = 0
$ws.Cells.Item(1501, 13).Value2 = 0
$ws.Cells.Item(1501, 14).Value2 = 39
$ws.Cells.Item(1501, 15).Value2 = 1
$ws.Cells.Item(1501, 16).Value2 = 0
$ws.Cells.Item(1501, 17).Value2 = 0

# Row 1502
$ws.Cells.Item(1502, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1502, 1).Value2 = 45565
$ws.Cells.Item(1502, 2).Value2 = 2964.199951171875
$ws.Cells.Item(1502, 3).Value2 = 2997.89990234375
$ws.Cells.Item(1502, 4).Value2 = 2817.050048828125
$ws.Cells.Item(1502, 5).Value2 = 2848.75
$ws.Cells.Item(1502, 6).Value2 = 2848.75
$ws.Cells.Item(1502, 7).Value2 = 6502277
$ws.Cells.Item(1502, 8).Value2 = 2024
$ws.Cells.Item(1502, 9).Value2 = 9
$ws.Cells.Item(1502, 10).Value2 = 30
$ws.Cells.Item(1502, 11).Value2 = 0
$ws.Cells.Item(1502, 12).Value2 = 0
$ws.Cells.Item(1502, 13).Value2 = 0
$ws.Cells.Item(1502, 14).Value2 = 40
$ws.Cells.Item(1502, 15).Value2 = 0
$ws.Cells.Item(1502, 16).Value2 = 0
$ws.Cells.Item(1502, 17).Value2 = 0

# Row 1503
$ws.Cells.Item(1503, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1503, 1).Value2 = 45572
$ws.Cells.Item(1503, 2).Value2 = 2845
$ws.Cells.Item(1503, 3).Value2 = 2886
$ws.Cells.Item(1503, 4).Value2 = 2733.199951171875
$ws.Cells.Item(1503, 5).Value2 = 2783.199951171875
$ws.Cells.Item(1503, 6).Value2 = 2783.199951171875
$ws.Cells.Item(1503, 7).Value2 = 7263130
$ws.Cells.Item(1503, 8).Value2 = 2024
$ws.Cells.Item(1503, 9).Value2 = 10
$ws.Cells.Item(1503, 10).Value2 = 7
$ws.Cells.Item(1503, 11).Value2 = 0
$ws.Cells.Item(1503, 12).Value2 = 0
$ws.Cells.Item(1503, 13).Value2 = 0
$ws.Cells.Item(1503, 14).Value2 = 41
$ws.Cells.Item(1503, 15).Value2 = 0
$ws.Cells.Item(1503, 16).Value2 = 0
$ws.Cells.Item(1503, 17).Value2 = 0

# Row 1504
$ws.Cells.Item(1504, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1504, 1).Value2 = 45579
$ws.Cells.Item(1504, 2).Value2 = 2791.60009765625
$ws.Cells.Item(1504, 3).Value2 = 2804.699951171875
$ws.Cells.Item(1504, 4).Value2 = 2693.449951171875
$ws.Cells.Item(1504, 5).Value2 = 2717.10009765625
$ws.Cells.Item(1504, 6).Value2 = 2717.10009765625
$ws.Cells.Item(1504, 7).Value2 = 5632992
$ws.Cells.Item(1504, 8).Value2 = 2024
$ws.Cells.Item(1504, 9).Value2 = 10
$ws.Cells.Item(1504, 10).Value2 = 14
$ws.Cells.Item(1504, 11).Value2 = 0
$ws.Cells.Item(1504, 12).Value2 = 0
$ws.Cells.Item(1504, 13).Value2 = 0
$ws.Cells.Item(1504, 14).Value2 = 42
$ws.Cells.Item(1504, 15).Value2 = 0
$ws.Cells.Item(1504, 16).Value2 = 0
$ws.Cells.Item(1504, 17).Value2 = 0

# Row 1505
$ws.Cells.Item(1505, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1505, 1).Value2 = 45586
$ws.Cells.Item(1505, 2).Value2 = 2715
$ws.Cells.Item(1505, 3).Value2 = 2738
$ws.Cells.Item(1505, 4).Value2 = 2452.60009765625
$ws.Cells.Item(1505, 5).Value2 = 2528.050048828125
$ws.Cells.Item(1505, 6).Value2 = 2528.050048828125
$ws.Cells.Item(1505, 7).Value2 = 15473993
$ws.Cells.Item(1505, 8).Value2 = 2024
$ws.Cells.Item(1505, 9).Value2 = 10
$ws.Cells.Item(1505, 10).Value2 = 21
$ws.Cells.Item(1505, 11).Value2 = 0
$ws.Cells.Item(1505, 12).Value2 = 0
$ws.Cells.Item(1505, 13).Value2 = 0
$ws.Cells.Item(1505, 14).Value2 = 43
$ws.Cells.Item(1505, 15).Value2 = 0
$ws.Cells.Item(1505, 16).Value2 = 0
$ws.Cells.Item(1505, 17).Value2 = 0

# Row 1506
$ws.Cells.Item(1506, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1506, 1).Value2 = 45593
$ws.Cells.Item(1506, 2).Value2 = 2544
$ws.Cells.Item(1506, 3).Value2 = 2589.60009765625
$ws.Cells.Item(1506, 4).Value2 = 2520.10009765625
$ws.Cells.Item(1506, 5).Value2 = 2528.25
$ws.Cells.Item(1506, 6).Value2 = 2528.25
$ws.Cells.Item(1506, 7).Value2 = 5650122
$ws.Cells.Item(1506, 8).Value2 = 2024
$ws.Cells.Item(1506, 9).Value2 = 10
$ws.Cells.Item(1506, 10).Value2 = 28
$ws.Cells.Item(1506, 11).Value2 = 0
$ws.Cells.Item(1506, 12).Value2 = 0
$ws.Cells.Item(1506, 13).Value2 = 0
$ws.Cells.Item(1506, 14).Value2 = 44
$ws.Cells.Item(1506, 15).Value2 = 0
$ws.Cells.Item(1506, 16).Value2 = 0
$ws.Cells.Item(1506, 17).Value2 = 0
